$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new row of data (row 2)
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "wondifraw"
$ws.Range("C2").Value = "nigussie"
$ws.Range("D2").Value = 45
$ws.Range("E2").Value = 911513683
$ws.Range("F2").Value = "aa"
$ws.Range("G2").Value = "ethiopia"

# Set column E width (bestFit) to match the new content (AutoFit-equivalent: ~10 chars wide)
$ws.Columns.Item(5).ColumnWidth = 9.14

# Update the active selection to A3 (as if the user pressed Enter/moved down after typing row 2)
$ws.Range("A3").Select()
